$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2040
$ws.Range("C2").Value = 149092
$ws.Range("D2").Value = 149092
$ws.Range("E2").Value = 44.98407253963344
$ws.Range("F2").Value = 0.2203999800807702
$ws.Range("G2").Value = [double]"3.694822225952521e-12"
$ws.Range("H2").Value = [double]"2.801243681460619e-10"
$ws.Range("I2").Value = 0.7560416917044475
$ws.Range("J2").Value = 0.5715990395953239
